$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Fix casing of existing label in A3
$ws.Range("A3").Value = "population"

# Add new row with density figures
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 8514.324289895707
